# effort_estimation.xlsx: add a "sigma" (standard deviation) column next to
# the existing weighted-mean "<T>" column, relabel the two headers, widen
# the columns that now hold longer text, and move the selection down one
# row to reflect the extra data row that was filled in (D2/D6 continued).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabelling (E1/F1) ------------------------------------------
# "<T>" -> "<T> (weighted mean)" and "sigmahoch2" -> "sigma"
$ws.Range("E1").Value = "<T> (weighted mean)"
$ws.Range("F1").Value = "sigma"

# --- New column F: sigma = (pessimistic - optimistic) / 6 ----------------
# F2 keeps its own (non-shared) formula, F3:F15 share one formula, matching
# how the existing E-column formulas are already laid out (E2 standalone,
# E3:E15 shared).
$ws.Range("F2").Formula = "=(D2-B2)/6"
$ws.Range("F3:F15").Formula = "=(D3-B3)/6"

# --- Column widths ----------------------------------------------------
# Give column D (pessimistic) and the now wider column E (weighted-mean
# header) more room.
$ws.Columns("D").ColumnWidth = 16.5
$ws.Columns("E").ColumnWidth = 25.333333333333332

# --- Selection moves from H15 to H16 --------------------------------------
[void]$ws.Range("H16").Select()
